# Versuch Verknüpfung Registrierung & Profil mit DB -> fehlerhaft
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: "Erstellen der DB" is now done, assigned to Nadine,
# but a foreign key is still missing.
$ws.Range("C4").Value = "erledigt"
$ws.Range("D4").Value = "Nadine"
$ws.Range("E4").Value = "Fremdschlüssel fehlt noch!"

# Update current selection to C5 (as left by the editor after the change)
$ws.Range("C5").Select()
